$d = $word.ActiveDocument

$d.Content.Find.Execute("854÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "398÷3=", 2)
$d.Content.Find.Execute("544÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "197÷5=", 2)
$d.Content.Find.Execute("213÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "715÷5=", 2)
$d.Content.Find.Execute("349÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "331÷5=", 2)
$d.Content.Find.Execute("162÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "980÷5=", 2)
$d.Content.Find.Execute("482÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "550÷2=", 2)
$d.Content.Find.Execute("360÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "801÷8=", 2)
$d.Content.Find.Execute("273÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "587÷3=", 2)
$d.Content.Find.Execute("624÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "976÷8=", 2)
$d.Content.Find.Execute("308÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "214÷3=", 2)
$d.Content.Find.Execute("763÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "540÷8=", 2)
$d.Content.Find.Execute("898÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "131÷3=", 2)
$d.Content.Find.Execute("157÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "824÷7=", 2)
$d.Content.Find.Execute("163÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "647÷8=", 2)
$d.Content.Find.Execute("417÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "229÷3=", 2)
$d.Content.Find.Execute("502÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "415÷9=", 2)
$d.Content.Find.Execute("956÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "285÷2=", 2)
$d.Content.Find.Execute("186÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "131÷5=", 2)
$d.Content.Find.Execute("996÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "714÷8=", 2)
$d.Content.Find.Execute("916÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "128÷6=", 2)
$d.Content.Find.Execute("447÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "877÷8=", 2)
$d.Content.Find.Execute("604÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "769÷3=", 2)
$d.Content.Find.Execute("756÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "992÷8=", 2)
$d.Content.Find.Execute("925÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "202÷5=", 2)
$d.Content.Find.Execute("467÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "383÷6=", 2)
